$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 615 (reverts an appended month row)
$ws.Rows("615:615").Delete()

# Apply updated cell values across rows 98-212 and 604-614
$ws.Range("H98").Value = 904865
$ws.Range("M99").Value = 8586174
$ws.Range("M100").Value = 8552547
$ws.Range("M102").Value = 8469973
$ws.Range("M104").Value = 8353737
$ws.Range("M105").Value = 8358614
$ws.Range("L105").Value = 82.49
$ws.Range("M109").Value = 8212235
$ws.Range("H132").Value = 936159
$ws.Range("M132").Value = 14067706
$ws.Range("M133").Value = 14470158
$ws.Range("L133").Value = 99.54000000000001
$ws.Range("M134").Value = 14605794
$ws.Range("M135").Value = 14691230
$ws.Range("L135").Value = 100.73
$ws.Range("M136").Value = 14867690
$ws.Range("M137").Value = 14862208
$ws.Range("M138").Value = 14671600
$ws.Range("M139").Value = 14467825
$ws.Range("D140").Value = 823660
$ws.Range("N140").Value = 8935496
$ws.Range("M140").Value = 14132605
$ws.Range("N141").Value = 8782011
$ws.Range("M141").Value = 13894901
$ws.Range("M142").Value = 13750808
$ws.Range("N142").Value = 8776279
$ws.Range("M143").Value = 13622677
$ws.Range("N143").Value = 8836721
$ws.Range("N144").Value = 9042529
$ws.Range("N145").Value = 9395149
$ws.Range("N146").Value = 9716013
$ws.Range("N147").Value = 9916844
$ws.Range("N148").Value = 10176477
$ws.Range("N149").Value = 10402058
$ws.Range("N150").Value = 10586458
$ws.Range("N151").Value = 10834371
$ws.Range("D186").Value = 590478
$ws.Range("H186").Value = 1208914
$ws.Range("N186").Value = 8211283
$ws.Range("M186").Value = 13889703
$ws.Range("M187").Value = 13965936
$ws.Range("N187").Value = 8231274
$ws.Range("N188").Value = 8311060
$ws.Range("M188").Value = 14177170
$ws.Range("N189").Value = 8315420
$ws.Range("M189").Value = 14260477
$ws.Range("L190").Value = 126.79
$ws.Range("N190").Value = 8371001
$ws.Range("M190").Value = 14422805
$ws.Range("D191").Value = 573529
$ws.Range("E191").Value = 207683
$ws.Range("H191").Value = 1073060
$ws.Range("O191").Value = 2551099
$ws.Range("N191").Value = 8368504
$ws.Range("M191").Value = 14527962
$ws.Range("N192").Value = 8290078
$ws.Range("O192").Value = 2547721
$ws.Range("M192").Value = 14532737
$ws.Range("O193").Value = 2572369
$ws.Range("M193").Value = 14761260
$ws.Range("N193").Value = 8365170
$ws.Range("N194").Value = 8388057
$ws.Range("O194").Value = 2570169
$ws.Range("M194").Value = 14843941
$ws.Range("M195").Value = 14888871
$ws.Range("N195").Value = 8325529
$ws.Range("O195").Value = 2571058
$ws.Range("L196").Value = 130.2
$ws.Range("N196").Value = 8320025
$ws.Range("M196").Value = 14967553
$ws.Range("O196").Value = 2563132
$ws.Range("M197").Value = 15089698
$ws.Range("O197").Value = 2566587
$ws.Range("N197").Value = 8341489
$ws.Range("N198").Value = 8330289
$ws.Range("O198").Value = 2567388
$ws.Range("M198").Value = 15141509
$ws.Range("M199").Value = 15309053
$ws.Range("N199").Value = 8356537
$ws.Range("O199").Value = 2587943
$ws.Range("N200").Value = 8350556
$ws.Range("O200").Value = 2599996
$ws.Range("M200").Value = 15440662
$ws.Range("E201").Value = 210443
$ws.Range("O201").Value = 2605573
$ws.Range("L201").Value = 133.49
$ws.Range("N201").Value = 8368299
$ws.Range("M201").Value = 15510437
$ws.Range("O202").Value = 2635485
$ws.Range("N202").Value = 8381856
$ws.Range("M202").Value = 15673747
$ws.Range("O203").Value = 2650975
$ws.Range("O204").Value = 2658797
$ws.Range("O205").Value = 2687346
$ws.Range("O206").Value = 2682400
$ws.Range("O207").Value = 2694751
$ws.Range("O208").Value = 2719334
$ws.Range("O209").Value = 2707776
$ws.Range("O210").Value = 2680341
$ws.Range("O211").Value = 2671978
$ws.Range("O212").Value = 2637236
$ws.Range("B604").Value = 11769240
$ws.Range("F604").Value = 14907682
$ws.Range("H604").Value = 3891902
$ws.Range("G604").Value = 10615357
$ws.Range("E604").Value = 164910
$ws.Range("M604").Value = 28616691
$ws.Range("O604").Value = 1773972
$ws.Range("H605").Value = 18338671
$ws.Range("E605").Value = 206161
$ws.Range("B605").Value = 17632868
$ws.Range("D605").Value = 12201277
$ws.Range("G605").Value = 55715641
$ws.Range("F605").Value = 66125880
$ws.Range("M605").Value = 44516864
$ws.Range("O605").Value = 1812735
$ws.Range("N605").Value = 18275692
$ws.Range("F606").Value = 87031320
$ws.Range("E606").Value = 263886
$ws.Range("D606").Value = 5760517
$ws.Range("H606").Value = 23754356
$ws.Range("B606").Value = 9283434
$ws.Range("G606").Value = 76904907
$ws.Range("O606").Value = 1925972
$ws.Range("N606").Value = 23685029
$ws.Range("M606").Value = 66265409
$ws.Range("B607").Value = 6650692
$ws.Range("F607").Value = 77521310
$ws.Range("H607").Value = 22138432
$ws.Range("G607").Value = 75658953
$ws.Range("E607").Value = 605364
$ws.Range("D607").Value = 3163674
$ws.Range("O607").Value = 2392651
$ws.Range("M607").Value = 86487991
$ws.Range("N607").Value = 26483303
$ws.Range("E608").Value = 950158
$ws.Range("F608").Value = 68602070
$ws.Range("H608").Value = 18662865
$ws.Range("G608").Value = 64868385
$ws.Range("B608").Value = 5952705
$ws.Range("D608").Value = 2163838
$ws.Range("O608").Value = 3184761
$ws.Range("M608").Value = 102836944
$ws.Range("N608").Value = 28174614
$ws.Range("B609").Value = 3754528
$ws.Range("F609").Value = 63219468
$ws.Range("G609").Value = 56553382
$ws.Range("E609").Value = 1066604
$ws.Range("H609").Value = 16358572
$ws.Range("N609").Value = 29141306
$ws.Range("M609").Value = 117204480
$ws.Range("O609").Value = 4115422
$ws.Range("D610").Value = 913791
$ws.Range("E610").Value = 3056341
$ws.Range("G610").Value = 46146973
$ws.Range("F610").Value = 49206838
$ws.Range("H610").Value = 13636573
$ws.Range("J610").Value = 109.42
$ws.Range("B610").Value = 3499689
$ws.Range("O610").Value = 7032526
$ws.Range("M610").Value = 128896991
$ws.Range("N610").Value = 29743907
$ws.Range("B611").Value = 3454273
$ws.Range("H611").Value = 8377585
$ws.Range("G611").Value = 28459956
$ws.Range("F611").Value = 32123268
$ws.Range("E611").Value = 2435546
$ws.Range("D611").Value = 762364
$ws.Range("O611").Value = 9326452
$ws.Range("N611").Value = 30141899
$ws.Range("M611").Value = 135364585
$ws.Range("B612").Value = 3436702
$ws.Range("F612").Value = 26461549
$ws.Range("G612").Value = 21945627
$ws.Range("E612").Value = 1329892
$ws.Range("H612").Value = 6701818
$ws.Range("M612").Value = 140240338
$ws.Range("N612").Value = 30517970
$ws.Range("O612").Value = 10528441
$ws.Range("H613").Value = 6405711
$ws.Range("B613").Value = 4226107
$ws.Range("F613").Value = 23050860
$ws.Range("G613").Value = 20380363
$ws.Range("E613").Value = 1042565
$ws.Range("D613").Value = 874957
$ws.Range("O613").Value = 11412095
$ws.Range("N613").Value = 30831490
$ws.Range("M613").Value = 143966702
$ws.Range("B614").Value = 4630041
$ws.Range("M614").Value = 146802057
$ws.Range("D614").Value = 942953
$ws.Range("E614").Value = 1602667.666666667
$ws.Range("H614").Value = 5782524
$ws.Range("C614").Value = 6468.066666666667
$ws.Range("O614").Value = 12103326
$ws.Range("F614").Value = 23069341
$ws.Range("N614").Value = 31037950
$ws.Range("I614").Value = 345.31
$ws.Range("G614").Value = 18036771
$ws.Range("J614").Value = 42.96
$ws.Range("L614").Value = 318.84
